$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Values (string literals written in the order that matches the
#     target shared-strings table: Clovis, helpMe, Quickoff,
#     Intro et Spec, Gilles, Design) ---

# Row 2: Clovis kicks off the "Quickoff" action for the helpMe project
$ws.Range("A2").Value = "Clovis"
$ws.Range("B2").Value = 42975
$ws.Range("H2").Value = "helpMe"
$ws.Range("C2").Value = "Quickoff"

# Row 3: Intro et Spec, assigned to Gilles, due 2017-09-04
$ws.Range("A3").Value = "Clovis"
$ws.Range("B3").Value = 42975
$ws.Range("C3").Value = "Intro et Spec"
$ws.Range("D3").Value = "Gilles"
$ws.Range("E3").Value = 42982
$ws.Range("H3").Value = "helpMe"

# Row 4: Design, assigned to Clovis, due 2017-09-04
$ws.Range("A4").Value = "Clovis"
$ws.Range("B4").Value = 42975
$ws.Range("C4").Value = "Design"
$ws.Range("D4").Value = "Clovis"
$ws.Range("E4").Value = 42982
$ws.Range("H4").Value = "helpMe"

# --- Formatting ---

# Date cells: apply the built-in short-date format to B2 first, then
# clone that exact style onto the other date cells via copy/paste-special
# so they all share a single cellXfs entry.
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)

# "state" cells left empty but flagged red, sharing one cellXfs entry.
$ws.Range("F3").Interior.Color = 255
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("G2").Select()
